$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.040469427188328
$ws.Cells.Item(2, 4).Value = 1.03405753044538
$ws.Cells.Item(2, 5).Value = 1.047758228971069
$ws.Cells.Item(2, 6).Value = 1.056043552281245
$ws.Cells.Item(2, 9).Value = 1.036670433895397
$ws.Cells.Item(2, 10).Value = 1.045555864413027
$ws.Cells.Item(2, 11).Value = 1.036858098106587
$ws.Cells.Item(2, 12).Value = 1.050519972059028
$ws.Cells.Item(2, 13).Value = 1.058782346459902
$ws.Cells.Item(2, 14).Value = 1.019020010373312

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.04162434526889
$ws.Cells.Item(3, 4).Value = 1.034567264428856
$ws.Cells.Item(3, 5).Value = 1.0488224980137
$ws.Cells.Item(3, 6).Value = 1.057278947973397
$ws.Cells.Item(3, 9).Value = 1.036876622630976
$ws.Cells.Item(3, 10).Value = 1.046355387686997
$ws.Cells.Item(3, 11).Value = 1.037177739477238
$ws.Cells.Item(3, 12).Value = 1.051395495823166
$ws.Cells.Item(3, 13).Value = 1.059830228762355
$ws.Cells.Item(3, 14).Value = 1.019291637777339

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.042371534510398
$ws.Cells.Item(4, 4).Value = 1.034897186743521
$ws.Cells.Item(4, 5).Value = 1.049511420318782
$ws.Cells.Item(4, 6).Value = 1.05807890629347
$ws.Cells.Item(4, 9).Value = 1.037008810324822
$ws.Cells.Item(4, 10).Value = 1.04687208065869
$ws.Cells.Item(4, 11).Value = 1.037383943545992
$ws.Cells.Item(4, 12).Value = 1.051961689464372
$ws.Cells.Item(4, 13).Value = 1.060508288731429
$ws.Cells.Item(4, 14).Value = 1.019467023181818

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.042685625450426
$ws.Cells.Item(5, 4).Value = 1.035035906536953
$ws.Cells.Item(5, 5).Value = 1.049801108321652
$ws.Cells.Item(5, 6).Value = 1.05841534738318
$ws.Cells.Item(5, 9).Value = 1.037064087387848
$ws.Cells.Item(5, 10).Value = 1.047089142871445
$ws.Cells.Item(5, 11).Value = 1.037470481452589
$ws.Cells.Item(5, 12).Value = 1.052199639186843
$ws.Cells.Item(5, 13).Value = 1.060793348392515
$ws.Cells.Item(5, 14).Value = 1.019540665350924

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.042738361131626
$ws.Cells.Item(6, 4).Value = 1.035059199365678
$ws.Cells.Item(6, 5).Value = 1.049849752046432
$ws.Cells.Item(6, 6).Value = 1.058471845480758
$ws.Cells.Item(6, 9).Value = 1.03707335136815
$ws.Cells.Item(6, 10).Value = 1.047125579469897
$ws.Cells.Item(6, 11).Value = 1.037485002717877
$ws.Cells.Item(6, 12).Value = 1.052239587431588
$ws.Cells.Item(6, 13).Value = 1.060841211379097
$ws.Cells.Item(6, 14).Value = 1.019553024928503

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.042375731515002
$ws.Cells.Item(7, 4).Value = 1.03489904024532
$ws.Cells.Item(7, 5).Value = 1.049515290887577
$ws.Cells.Item(7, 6).Value = 1.058083401289837
$ws.Cells.Item(7, 9).Value = 1.037009550097554
$ws.Cells.Item(7, 10).Value = 1.046874981664925
$ws.Cells.Item(7, 11).Value = 1.037385100461001
$ws.Cells.Item(7, 12).Value = 1.051964869266247
$ws.Cells.Item(7, 13).Value = 1.060512097698745
$ws.Cells.Item(7, 14).Value = 1.019468007545049

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.040859762790886
$ws.Cells.Item(8, 4).Value = 1.034229778020917
$ws.Cells.Item(8, 5).Value = 1.048117848356079
$ws.Cells.Item(8, 6).Value = 1.056460941662834
$ws.Cells.Item(8, 9).Value = 1.036740370988414
$ws.Cells.Item(8, 10).Value = 1.045826202412947
$ws.Cells.Item(8, 11).Value = 1.036966251515923
$ws.Cells.Item(8, 12).Value = 1.050815927867176
$ws.Cells.Item(8, 13).Value = 1.059136481617488
$ws.Cells.Item(8, 14).Value = 1.019111885967245

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.038187454497957
$ws.Cells.Item(9, 4).Value = 1.03305119596516
$ws.Cells.Item(9, 5).Value = 1.045657407262588
$ws.Cells.Item(9, 6).Value = 1.053606315657373
$ws.Cells.Item(9, 9).Value = 1.036256626609897
$ws.Cells.Item(9, 10).Value = 1.043973098680313
$ws.Cells.Item(9, 11).Value = 1.03622342322234
$ws.Cells.Item(9, 12).Value = 1.048788799861154
$ws.Cells.Item(9, 13).Value = 1.056712509142609
$ws.Cells.Item(9, 14).Value = 1.018481472848548

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.036405173670799
$ws.Cells.Item(10, 4).Value = 1.03226604324953
$ws.Cells.Item(10, 5).Value = 1.044018434816761
$ws.Cells.Item(10, 6).Value = 1.051706092326859
$ws.Cells.Item(10, 9).Value = 1.035927805657276
$ws.Cells.Item(10, 10).Value = 1.042734280807751
$ws.Cells.Item(10, 11).Value = 1.035725036074683
$ws.Cells.Item(10, 12).Value = 1.0474356269876
$ws.Cells.Item(10, 13).Value = 1.055096497726269
$ws.Cells.Item(10, 14).Value = 1.018059253427595

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.035633230815336
$ws.Cells.Item(11, 4).Value = 1.031926212602915
$ws.Cells.Item(11, 5).Value = 1.043309044467686
$ws.Cells.Item(11, 6).Value = 1.050883937300936
$ws.Cells.Item(11, 9).Value = 1.035783923395868
$ws.Cells.Item(11, 10).Value = 1.042197039325576
$ws.Cells.Item(11, 11).Value = 1.03550848435843
$ws.Cells.Item(11, 12).Value = 1.046849261998886
$ws.Cells.Item(11, 13).Value = 1.054396726652747
$ws.Cells.Item(11, 14).Value = 1.01787596484268

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.03534646450689
$ws.Cells.Item(12, 4).Value = 1.03180000738221
$ws.Cells.Item(12, 5).Value = 1.043045588689996
$ws.Cells.Item(12, 6).Value = 1.05057864914139
$ws.Cells.Item(12, 9).Value = 1.035730253688878
$ws.Cells.Item(12, 10).Value = 1.041997359025582
$ws.Cells.Item(12, 11).Value = 1.035427935729794
$ws.Cells.Item(12, 12).Value = 1.046631394034383
$ws.Cells.Item(12, 13).Value = 1.054136795154715
$ws.Cells.Item(12, 14).Value = 1.017807813274197

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.035407978346642
$ws.Cells.Item(13, 4).Value = 1.031827077771073
$ws.Cells.Item(13, 5).Value = 1.043102098901796
$ws.Cells.Item(13, 6).Value = 1.050644130130011
$ws.Cells.Item(13, 9).Value = 1.035741776223989
$ws.Cells.Item(13, 10).Value = 1.042040196794227
$ws.Cells.Item(13, 11).Value = 1.035445218729516
$ws.Cells.Item(13, 12).Value = 1.046678130449972
$ws.Cells.Item(13, 13).Value = 1.054192551606023
$ws.Cells.Item(13, 14).Value = 1.017822435189049

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.035609527293328
$ws.Cells.Item(14, 4).Value = 1.03191577996406
$ws.Cells.Item(14, 5).Value = 1.043287266251799
$ws.Cells.Item(14, 6).Value = 1.050858700119742
$ws.Cells.Item(14, 9).Value = 1.035779491639698
$ws.Cells.Item(14, 10).Value = 1.042180536242796
$ws.Cells.Item(14, 11).Value = 1.035501828458903
$ws.Cells.Item(14, 12).Value = 1.046831254314906
$ws.Cells.Item(14, 13).Value = 1.054375240750622
$ws.Cells.Item(14, 14).Value = 1.017870332842451

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.035733703909658
$ws.Cells.Item(15, 4).Value = 1.031970435395678
$ws.Cells.Item(15, 5).Value = 1.043401359685701
$ws.Cells.Item(15, 6).Value = 1.050990916558329
$ws.Cells.Item(15, 9).Value = 1.03580269948124
$ws.Cells.Item(15, 10).Value = 1.042266987431763
$ws.Cells.Item(15, 11).Value = 1.035536692798366
$ws.Cells.Item(15, 12).Value = 1.04692559022053
$ws.Cells.Item(15, 13).Value = 1.054487800814284
$ws.Cells.Item(15, 14).Value = 1.017899834882319

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.036456400535751
$ws.Cells.Item(16, 4).Value = 1.032288599855766
$ws.Cells.Item(16, 5).Value = 1.044065520857723
$ws.Cells.Item(16, 6).Value = 1.051760669706792
$ws.Cells.Item(16, 9).Value = 1.035937323034703
$ws.Cells.Item(16, 10).Value = 1.042769918291673
$ws.Cells.Item(16, 11).Value = 1.035739392209139
$ws.Cells.Item(16, 12).Value = 1.047474532912291
$ws.Cells.Item(16, 13).Value = 1.055142938543502
$ws.Cells.Item(16, 14).Value = 1.018071407877905

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.036909673292835
$ws.Cells.Item(17, 4).Value = 1.032488215812649
$ws.Cells.Item(17, 5).Value = 1.044482210020287
$ws.Cells.Item(17, 6).Value = 1.052243689639146
$ws.Cells.Item(17, 9).Value = 1.036021367090834
$ws.Cells.Item(17, 10).Value = 1.043085171962893
$ws.Cells.Item(17, 11).Value = 1.035866340623725
$ws.Cells.Item(17, 12).Value = 1.047818753765545
$ws.Cells.Item(17, 13).Value = 1.055553880906366
$ws.Cells.Item(17, 14).Value = 1.018178906508793

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.037174040108892
$ws.Cells.Item(18, 4).Value = 1.032604662401966
$ws.Cells.Item(18, 5).Value = 1.044725286342274
$ws.Cells.Item(18, 6).Value = 1.052525490214346
$ws.Cells.Item(18, 9).Value = 1.036070243777435
$ws.Cells.Item(18, 10).Value = 1.043268974540722
$ws.Cells.Item(18, 11).Value = 1.035940315397017
$ws.Cells.Item(18, 12).Value = 1.048019490328888
$ws.Cells.Item(18, 13).Value = 1.055793574067048
$ws.Cells.Item(18, 14).Value = 1.018241563809373

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.037264179105708
$ws.Cells.Item(19, 4).Value = 1.032644370016279
$ws.Cells.Item(19, 5).Value = 1.04480817396554
$ws.Cells.Item(19, 6).Value = 1.052621587695971
$ws.Cells.Item(19, 9).Value = 1.03608688490499
$ws.Cells.Item(19, 10).Value = 1.043331633003955
$ws.Cells.Item(19, 11).Value = 1.035965526639643
$ws.Cells.Item(19, 12).Value = 1.048087929276551
$ws.Cells.Item(19, 13).Value = 1.055875302842826
$ws.Cells.Item(19, 14).Value = 1.018262920731353

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.036861043446859
$ws.Cells.Item(20, 4).Value = 1.032466797455525
$ws.Cells.Item(20, 5).Value = 1.0444375002673
$ws.Cells.Item(20, 6).Value = 1.052191859642793
$ws.Cells.Item(20, 9).Value = 1.036012364935035
$ws.Cells.Item(20, 10).Value = 1.043051356462717
$ws.Cells.Item(20, 11).Value = 1.035852727713001
$ws.Cells.Item(20, 12).Value = 1.04778182643913
$ws.Cells.Item(20, 13).Value = 1.055509790954195
$ws.Cells.Item(20, 14).Value = 1.018167377561618

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.035550177059699
$ws.Cells.Item(21, 4).Value = 1.031889658734792
$ws.Cells.Item(21, 5).Value = 1.043232737892982
$ws.Cells.Item(21, 6).Value = 1.050795511944752
$ws.Cells.Item(21, 9).Value = 1.035768391612565
$ws.Cells.Item(21, 10).Value = 1.042139213228118
$ws.Cells.Item(21, 11).Value = 1.035485161383314
$ws.Cells.Item(21, 12).Value = 1.046786164980118
$ws.Cells.Item(21, 13).Value = 1.054321443508965
$ws.Cells.Item(21, 14).Value = 1.017856230112029

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.034725795215973
$ws.Cells.Item(22, 4).Value = 1.031526922409687
$ws.Cells.Item(22, 5).Value = 1.042475506469878
$ws.Cells.Item(22, 6).Value = 1.049918132190795
$ws.Cells.Item(22, 9).Value = 1.035613691650738
$ws.Cells.Item(22, 10).Value = 1.041564989525074
$ws.Cells.Item(22, 11).Value = 1.035253411954698
$ws.Cells.Item(22, 12).Value = 1.046159771347718
$ws.Cells.Item(22, 13).Value = 1.053574250656398
$ws.Cells.Item(22, 14).Value = 1.01766019410674

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.035162834058212
$ws.Cells.Item(23, 4).Value = 1.031719202751471
$ws.Cells.Item(23, 5).Value = 1.042876905770671
$ws.Cells.Item(23, 6).Value = 1.050383195257936
$ws.Cells.Item(23, 9).Value = 1.035695824648208
$ws.Cells.Item(23, 10).Value = 1.041869465219425
$ws.Cells.Item(23, 11).Value = 1.035376327762645
$ws.Cells.Item(23, 12).Value = 1.046491870973177
$ws.Cells.Item(23, 13).Value = 1.053970355138613
$ws.Cells.Item(23, 14).Value = 1.017764154982168

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.036883017240804
$ws.Cells.Item(24, 4).Value = 1.03247647544608
$ws.Cells.Item(24, 5).Value = 1.044457702591508
$ws.Cells.Item(24, 6).Value = 1.05221527919137
$ws.Cells.Item(24, 9).Value = 1.036016433069066
$ws.Cells.Item(24, 10).Value = 1.043066636477176
$ws.Cells.Item(24, 11).Value = 1.035858879024159
$ws.Cells.Item(24, 12).Value = 1.047798512436748
$ws.Cells.Item(24, 13).Value = 1.055529713313296
$ws.Cells.Item(24, 14).Value = 1.018172587135024

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.038878434819709
$ws.Cells.Item(25, 4).Value = 1.03335579245298
$ws.Cells.Item(25, 5).Value = 1.046293254081651
$ws.Cells.Item(25, 6).Value = 1.054343795741881
$ws.Cells.Item(25, 9).Value = 1.036382801360267
$ws.Cells.Item(25, 10).Value = 1.044452769478076
$ws.Cells.Item(25, 11).Value = 1.036416023050705
$ws.Cells.Item(25, 12).Value = 1.049313167138548
$ws.Cells.Item(25, 13).Value = 1.057339165277158
$ws.Cells.Item(25, 14).Value = 1.018644791897022
